$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.014.11"
$ws.Range("E2").Value = "  -1.86%  "
$ws.Range("D3").Value = "1.554.59"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'1.000"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "'286.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").Value = "'0.3789"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.65%  "
$ws.Range("D8").Value = "'0.3242"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.28%  "
$ws.Range("D9").Value = "'41.22"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -12.92%  "
$ws.Range("E10").Value = "  -3.67%  "
$ws.Range("D11").Value = "'0.07310"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.47%  "
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "'19.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.64%  "
$ws.Range("D14").Value = "'5.727"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.41%  "
$ws.Range("D15").Value = "'6.790"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.42%  "
$ws.Range("D16").Value = "1.566.36"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").Value = "'0.00001088"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.19%  "
$ws.Range("E18").Value = "  -1.66%  "
$ws.Range("D19").Value = "'85.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.11%  "
$ws.Range("D20").Value = "'6.418"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "'0.9987"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Value = "'15.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.95%  "
$ws.Range("D23").Value = "'11.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.01%  "
$ws.Range("D24").Value = "22.017.37"
$ws.Range("E24").Value = "  -1.82%  "
$ws.Range("E25").Value = "  -3.59%  "
$ws.Range("D26").Value = "'2.527"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.38%  "
$ws.Range("D27").Value = "'147.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.87%  "
$ws.Range("D28").Value = "'18.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.69%  "
$ws.Range("D29").Value = "'4.859"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("D30").Value = "1.741.32"
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").Value = "'120.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.48%  "
$ws.Range("D32").Value = "'1.108"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.61%  "
$ws.Range("D33").Value = "'5.932"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.60%  "
$ws.Range("D34").Value = "'1.648"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -16.54%  "
$ws.Range("D35").Value = "'0.08148"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.92%  "
$ws.Range("D36").Value = "'9.215"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.72%  "
$ws.Range("D37").Value = "'5.242"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.51%  "
$ws.Range("D38").Value = "'0.06188"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.85%  "
$ws.Range("D39").Value = "'0.02281"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.14%  "
$ws.Range("E40").Value = "  -4.24%  "
$ws.Range("D41").Value = "'1.220"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.83%  "
$ws.Range("D42").Value = "'10.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.35%  "
$ws.Range("D43").Value = "'0.9995"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").Value = "'0.5933"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.51%  "
$ws.Range("D45").Value = "'13.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.82%  "
$ws.Range("D46").Value = "'3.720"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.45%  "
$ws.Range("D47").Value = "'0.5732"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.83%  "
$ws.Range("D48").Value = "'1.934"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.04%  "
$ws.Range("D49").Value = "'119.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.89%  "
$ws.Range("E50").Value = "  -3.17%  "
$ws.Range("D51").Value = "'0.06879"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.31%  "
